# Update SOHP Excel template
# - add "Interviewee Date" / "Interviewer Date" columns to the name groups
# - merge the separate Identifier + "Interview Number" note into a single
#   mods:identifier element carrying a displayLabel
# - append a genre element to the typeOfResource string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Interviewee: insert 3 columns after the "Interviewee" label's
#        closing tag (old column F) for the new namePart/date pair ---
$ws.Columns("G:I").Insert()

# F1 keeps the opening close of </mods:namePart>, but the role-closing
# remainder moves out to the new I1 cell.
$ws.Range("F1").Value = "</mods:namePart>"
$ws.Range("I1").Value = "</mods:namePart><mods:role><mods:roleTerm authority=""marcrelator"" type=""text"">Interviewee</mods:roleTerm></mods:role></mods:name>"

$ws.Range("G1").Value = "<mods:namePart type=""date"">"
$ws.Range("H1").Value = "Interviewee Date"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Font.Size = 12

# --- 2. Interviewer: same treatment. After the previous insert the old
#        interviewer-role-closing cell (formerly I1) now lives at L1 ---
$ws.Columns("M:O").Insert()

$ws.Range("L1").Value = "</mods:namePart>"
$ws.Range("O1").Value = "</mods:namePart><mods:role><mods:roleTerm authority=""marcrelator"" type=""text"">Interviewer</mods:roleTerm></mods:role></mods:name>"

$ws.Range("M1").Value = "<mods:namePart type=""date"">"
$ws.Range("N1").Value = "Interviewer Date"
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").Font.Size = 12

# --- 3. Identifier / Interview Number: merge the old
#        <mods:identifier type="local"> + "Identifier" + "Interview
#        Number" note trio into one identifier carrying a displayLabel,
#        then drop the now-unused note columns ---
# After the two insertions above, the old P:W range shifted to V:AC.
$ws.Range("V1").Value = "<mods:identifier displayLabel=""Interview Number"" type=""local"">"
$ws.Range("W1").Value = "Interview Number"
# X1 already holds "</mods:identifier>" (shifted from old R1) - keep it.

# Remove the old "<mods:note displayLabel=...">, "Interview Number",
# "</mods:note>" trio (shifted to Y1:AA1); typeOfResource/mods-close
# shift back left into Y1/Z1.
$ws.Range("Y1:AA1").EntireColumn.Delete()

# --- 4. typeOfResource gains a trailing genre element ---
$ws.Range("Y1").Value = "<mods:typeOfResource>sound recording-nonmusical</mods:typeOfResource><mods:genre authority=""lcgft"">Oral histories</mods:genre>"

# --- 5. Sheet view: scrolled right with Z1 selected, matching the
#        widened header row ---
$ws.Range("Z1").Select()
